# Apply the two changes described by the commit:
#   1. Slide 5's table gets a new table style (tableStyleId).
#   2. The presentation's theme colors are swapped back to the previous
#      ("Office Theme") palette that used to live in ppt/theme/theme1.xml,
#      so the slide master's theme (ppt/theme/theme2.xml) now carries the
#      plain Office colors instead of the Integral/Red Violet ones.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 --------------------------------------------
$s = $p.Slides.Item(5)
$tbl = $s.Shapes.Item(2).Table
$tbl.ApplyStyle("{76257238-3120-46B2-A4E9-5530A136BCBB}")

# --- 2. Swap the theme color scheme back to the Office palette -----------
# ThemeColorScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# RGB values use the OLE (BGR-packed) integer form PowerPoint expects.
$officeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
